$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the prior row (row 24) onto the new row 25 first.
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)

# Now set the new row of timelog data. Set B (work) before A (date) so the
# shared-string insertion order matches the target (work text gets the
# lower index, 48, then the date string gets 49).
$ws.Range("B25").Value = "Working on graph over time for country, addinggoverment to said graph"
$ws.Range("A25").Value = "3/9, 4 hrs"

$ws.Rows("25").RowHeight = $ws.Rows("24").RowHeight

# Update view/selection to match new state
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("A26").Select()
